$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1145.8462
$ws.Range("I2").Value = 674.875
$ws.Range("J2").Value = 1899.4
$ws.Range("K2").Value = 674.875
$ws.Range("L2").Value = 1899.4
$ws.Range("M2").Value = -561.875
$ws.Range("N2").Value = -2125.4

$ws.Range("H6").Value = 195.5
$ws.Range("I6").Value = 195.5
$ws.Range("K6").Value = 586.5
$ws.Range("M6").Value = -474.5

$ws.Range("H33").Value = 341.95456
$ws.Range("I33").Value = 240.61111
$ws.Range("K33").Value = 240.61111
$ws.Range("M33").Value = -11.61111

$ws.Range("H51").Value = 7603.5713
$ws.Range("I51").Value = 20225
$ws.Range("K51").Value = 20225
$ws.Range("M51").Value = -19741

$ws.Range("H58").Value = 4999.4
$ws.Range("I58").Value = 498.5
$ws.Range("K58").Value = 1495.5
$ws.Range("M58").Value = -1345.5

$ws.Range("H62").Value = 4138.5
$ws.Range("I62").Value = 1482.1111
$ws.Range("K62").Value = 1482.1111
$ws.Range("M62").Value = -858.1111000000001

$ws.Range("H65").Value = 4138.5
$ws.Range("I65").Value = 1482.1111
$ws.Range("K65").Value = 7410.5555
$ws.Range("M65").Value = -4290.5555

$ws.Range("H86").Value = 3980.8096
$ws.Range("I86").Value = 1727.0834
$ws.Range("J86").Value = 6985.778
$ws.Range("K86").Value = 1727.0834
$ws.Range("L86").Value = 6985.778
$ws.Range("M86").Value = -604.0834
$ws.Range("N86").Value = -9231.778

$ws.Range("H89").Value = 3980.8096
$ws.Range("I89").Value = 1727.0834
$ws.Range("J89").Value = 6985.778
$ws.Range("K89").Value = 8635.416999999999
$ws.Range("L89").Value = 34928.89
$ws.Range("M89").Value = -3019.416999999999
$ws.Range("N89").Value = -46160.89

$ws.Range("H112").Value = 3060.3333
$ws.Range("I112").Value = 1397
$ws.Range("J112").Value = 3393
$ws.Range("K112").Value = 4191
$ws.Range("L112").Value = 10179
$ws.Range("M112").Value = -3083
$ws.Range("N112").Value = -12395

$ws.Range("H125").Value = 4069
$ws.Range("J125").Value = 6215.1665
$ws.Range("L125").Value = 55936.4985
$ws.Range("N125").Value = -60856.4985

$ws.Range("H138").Value = 4272.76
$ws.Range("I138").Value = 1697.8636
$ws.Range("K138").Value = 5093.5908
$ws.Range("M138").Value = 46.40920000000006


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4301.4
$ws.Range("I32").Value = 4319.378
$ws.Range("K32").Value = 4319.378
$ws.Range("M32").Value = -4032.378

$ws.Range("H45").Value = 1548.375
$ws.Range("I45").Value = 1484
$ws.Range("K45").Value = 1484
$ws.Range("M45").Value = -1107

$ws.Range("H61").Value = 70003336
$ws.Range("I61").Value = 100005000
$ws.Range("K61").Value = 100005000
$ws.Range("M61").Value = -100004788

$ws.Range("H102").Value = 4199.9287
$ws.Range("I102").Value = 3384.5
$ws.Range("K102").Value = 3384.5
$ws.Range("M102").Value = -1762.5

$ws.Range("H109").Value = 20599.25
$ws.Range("J109").Value = 20599.25
$ws.Range("L109").Value = 20599.25
$ws.Range("N109").Value = -23373.25

$ws.Range("H110").Value = 8216.444
$ws.Range("I110").Value = 7706.857
$ws.Range("K110").Value = 7706.857
$ws.Range("M110").Value = -5661.857

$ws.Range("H132").Value = 1727215.9
$ws.Range("I132").Value = 2961.9387
$ws.Range("J132").Value = 11114821
$ws.Range("K132").Value = 8885.8161
$ws.Range("L132").Value = 33344463
$ws.Range("M132").Value = -6355.8161
$ws.Range("N132").Value = -33349523

$ws.Range("H136").Value = 70003336
$ws.Range("I136").Value = 100005000
$ws.Range("K136").Value = 300015000
$ws.Range("M136").Value = -300012450


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 9097937
$ws.Range("I16").Value = 20004720
$ws.Range("J16").Value = 8951.666999999999
$ws.Range("K16").Value = 20004720
$ws.Range("L16").Value = 8951.666999999999
$ws.Range("M16").Value = -20004433
$ws.Range("N16").Value = -9525.666999999999

$ws.Range("H31").Value = 77711740
$ws.Range("I31").Value = 100004770
$ws.Range("K31").Value = 100004770
$ws.Range("M31").Value = -100004475

$ws.Range("H34").Value = 77711740
$ws.Range("I34").Value = 100004770
$ws.Range("K34").Value = 100004770
$ws.Range("M34").Value = -100004568

$ws.Range("H58").Value = 2372.1904
$ws.Range("I58").Value = 2201.8948
$ws.Range("K58").Value = 2201.8948
$ws.Range("M58").Value = -1998.8948

$ws.Range("H113").Value = 9097937
$ws.Range("I113").Value = 20004720
$ws.Range("J113").Value = 8951.666999999999
$ws.Range("K113").Value = 20004720
$ws.Range("L113").Value = 8951.666999999999
$ws.Range("M113").Value = -20002550
$ws.Range("N113").Value = -13291.667

$ws.Range("H122").Value = 3391.7646
$ws.Range("I122").Value = 3187.3333
$ws.Range("K122").Value = 9561.999899999999
$ws.Range("M122").Value = -7111.999899999999

$ws.Range("H132").Value = 2757.4167
$ws.Range("I132").Value = 2747.1177
$ws.Range("K132").Value = 8241.3531
$ws.Range("M132").Value = -5711.3531

$ws.Range("H134").Value = 2049.5715
$ws.Range("I134").Value = 2056.1667
$ws.Range("K134").Value = 6168.500100000001
$ws.Range("M134").Value = -3633.500100000001

$ws.Range("H136").Value = 2372.1904
$ws.Range("I136").Value = 2201.8948
$ws.Range("K136").Value = 6605.6844
$ws.Range("M136").Value = -4055.6844


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 27778278
$ws.Range("I7").Value = 1001
$ws.Range("J7").Value = 55555556
$ws.Range("K7").Value = 3003
$ws.Range("L7").Value = 166666668
$ws.Range("M7").Value = -2891
$ws.Range("N7").Value = -166666892

$ws.Range("H38").Value = 126
$ws.Range("I38").Value = 7.5
$ws.Range("K38").Value = 22.5
$ws.Range("M38").Value = 324.5

$ws.Range("H80").Value = 111114180
$ws.Range("J80").Value = 125003850
$ws.Range("L80").Value = 375011550
$ws.Range("N80").Value = -375013422

$ws.Range("H83").Value = 111114180
$ws.Range("J83").Value = 125003850
$ws.Range("L83").Value = 1125034650
$ws.Range("N83").Value = -1125044010

$ws.Range("H92").Value = 37037384
$ws.Range("I92").Value = 166666910
$ws.Range("J92").Value = 375.85715
$ws.Range("K92").Value = 500000730
$ws.Range("L92").Value = 1127.57145
$ws.Range("M92").Value = -499999482
$ws.Range("N92").Value = -3623.57145

$ws.Range("H107").Value = 4794047.5
$ws.Range("J107").Value = 7005211
$ws.Range("L107").Value = 21015633
$ws.Range("N107").Value = -21019473

$ws.Range("H137").Value = 5006.3
$ws.Range("I137").Value = 2887.4666
$ws.Range("J137").Value = 11362.8
$ws.Range("K137").Value = 8662.399800000001
$ws.Range("L137").Value = 34088.39999999999
$ws.Range("M137").Value = -3562.399800000001
$ws.Range("N137").Value = -44288.39999999999


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2932.889
$ws.Range("I102").Value = 2943
$ws.Range("K102").Value = 2943
$ws.Range("M102").Value = -1321

$ws.Range("H122").Value = 4334.615
$ws.Range("I122").Value = 4778.857
$ws.Range("K122").Value = 14336.571
$ws.Range("M122").Value = -11886.571

$ws.Range("H132").Value = 3574464.5
$ws.Range("J132").Value = 12503211
$ws.Range("L132").Value = 37509633
$ws.Range("N132").Value = -37514693


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4311.25
$ws.Range("I40").Value = 4212.857
$ws.Range("K40").Value = 4212.857
$ws.Range("M40").Value = -4076.857

$ws.Range("H46").Value = 8332.666999999999
$ws.Range("I46").Value = 9998
$ws.Range("K46").Value = 9998
$ws.Range("M46").Value = -9810

$ws.Range("H55").Value = 1425.875
$ws.Range("I55").Value = 1908
$ws.Range("K55").Value = 1908
$ws.Range("M55").Value = -1735

$ws.Range("H101").Value = 29749
$ws.Range("J101").Value = 29749
$ws.Range("L101").Value = 29749
$ws.Range("N101").Value = -36239

$ws.Range("H122").Value = 3548.5112
$ws.Range("I122").Value = 3402.375
$ws.Range("J122").Value = 4717.6
$ws.Range("K122").Value = 10207.125
$ws.Range("L122").Value = 14152.8
$ws.Range("M122").Value = -7757.125
$ws.Range("N122").Value = -19052.8

$ws.Range("H136").Value = 4006.4666
$ws.Range("I136").Value = 3354.6667
$ws.Range("K136").Value = 10064.0001
$ws.Range("M136").Value = -7514.000100000001


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 11836.9
$ws.Range("J96").Value = 11399.333
$ws.Range("L96").Value = 11399.333
$ws.Range("N96").Value = -14145.333

$ws.Range("H109").Value = 88073.336
$ws.Range("J109").Value = 99610.5
$ws.Range("L109").Value = 99610.5
$ws.Range("N109").Value = -102384.5

$ws.Range("H122").Value = 2970.0908
$ws.Range("I122").Value = 2505.6924
$ws.Range("K122").Value = 7517.0772
$ws.Range("M122").Value = -5067.0772

$ws.Range("H132").Value = 478692.34
$ws.Range("I132").Value = 2633.6316
$ws.Range("K132").Value = 7900.8948
$ws.Range("M132").Value = -5370.8948
